$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the read-control scores for the first three students
$ws.Range("C15").Value = 10
$ws.Range("D15").Value = 9
$ws.Range("E15").Value = 9
$ws.Range("F15").Value = 10

$ws.Range("C16").Value = 10
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = 9
$ws.Range("F16").Value = 10

$ws.Range("C17").Value = 9
$ws.Range("D17").Value = 9
$ws.Range("E17").Value = 10

# Update the view: scroll back to top and move selection to G17
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
$ws.Range("G17").Select()
